$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for Memory Usage column
$ws.Range("D1").Value = "Memory Usage (bytes)"

# Update Run Time (ms) values in column C
$ws.Range("C2").Value = 17.05098152160645
$ws.Range("C3").Value = 16.79897308349609
$ws.Range("C4").Value = 17.92216300964355
$ws.Range("C5").Value = 18.52703094482422
$ws.Range("C6").Value = 17.30990409851074
